$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 3..42 down from rows 2..41 (working bottom-up so sources
# aren't clobbered before they're read). This reproduces a new record
# being inserted as the new row 2, with existing rows 2-41 pushed down
# to rows 3-42. Rows 43+ are left untouched.
for ($i = 42; $i -ge 3; $i--) {
    $srcRow = $i - 1
    $ws.Range("A$i").Value = $ws.Range("A$srcRow").Value2
    $ws.Range("B$i").Value = $ws.Range("B$srcRow").Value2
    $ws.Range("C$i").Value = $ws.Range("C$srcRow").Value2
    $ws.Range("D$i").Value = $ws.Range("D$srcRow").Value2
}

# New record for row 2.
$ws.Range("A2").Value = "Springvale Shopping Centre, 46-58 Buckingham Avenue, Springvale VIC 3171, Australia"
$ws.Range("B2").Value = -37.9506608
$ws.Range("C2").Value = 145.1505924
$ws.Range("D2").Value = "Greater Dandenong (C)"
